# Rename ConstValue.getLiteralIntValue() -> getIntValue() on the two
# slides that reference it ("Code Generation for ConstValue" and
# "Method emit() for Class ConstValue").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 28 : "Class ConstValue has a method getLiteralIntValue() ..."
# ---------------------------------------------------------------
$slide28 = $p.Slides.Item(28)
$body28  = $slide28.Shapes.Item(2).TextFrame.TextRange

for ($i = 1; $i -le $body28.Paragraphs().Count; $i++) {
    $para = $body28.Paragraphs($i)
    $pos = $para.Text.IndexOf("getLiteralIntValue")
    if ($pos -ge 0) {
        $run = $para.Characters($pos + 1, "getLiteralIntValue".Length)
        $run.Text = "getIntValue"
    }
}

# ---------------------------------------------------------------
# Slide 29 : "emit(""LDCINT "" + getLiteralIntValue());"
#            "emit(""LDCB "" + getLiteralIntValue());"
# ---------------------------------------------------------------
$slide29 = $p.Slides.Item(29)
$body29  = $slide29.Shapes.Item(2).TextFrame.TextRange

for ($i = 1; $i -le $body29.Paragraphs().Count; $i++) {
    $para = $body29.Paragraphs($i)
    $pos = $para.Text.IndexOf("getLiteralIntValue")
    if ($pos -ge 0) {
        $run = $para.Characters($pos + 1, "getLiteralIntValue".Length)
        $run.Text = "getIntValue"
    }
}
